$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.416.04"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.849.29"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'240.84"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'0.6328"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.322.13"
$ws.Range("E8").Value = "  +79.70%  "
$ws.Range("D9").Value = "'0.07566"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'0.2972"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'24.63"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "3.660.80"
$ws.Range("E12").Value = "  +75.25%  "
$ws.Range("D13").Value = "'0.07712"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'4.986"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "'0.6843"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "'83.03"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "'0.000009956"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").Value = "'6.170"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "29.433.89"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'231.75"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'7.584"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'154.95"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").Value = "'0.1385"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "'8.415"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "'1.468"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'0.05787"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "'1.260"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "'4.129"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "3.615.97"
$ws.Range("E33").Value = "  +80.46%  "
$ws.Range("D34").Value = "'4.022"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'1.867"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "'0.7167"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'2.590"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "1.251.08"
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.01805"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "'0.9016"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "'6.087"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "'0.9992"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'101.69"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "'66.95"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "'7.189"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'9.148"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "'0.4012"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("E51").Value = "  +0.20%  "
